$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 242876.8
$ws.Range("J19").Value = 980.7
$ws.Range("L19").Value = 980.7
$ws.Range("N19").Value = -1330.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1265.5
$ws.Range("I28").Value = 1367.2727
$ws.Range("J28").Value = 892.3333
$ws.Range("K28").Value = 1367.2727
$ws.Range("L28").Value = 892.3333
$ws.Range("M28").Value = -882.2727
$ws.Range("N28").Value = -1862.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 106380.5
$ws.Range("I62").Value = 203959.8
$ws.Range("J62").Value = 8801.200000000001
$ws.Range("K62").Value = 203959.8
$ws.Range("L62").Value = 8801.200000000001
$ws.Range("M62").Value = -203335.8
$ws.Range("N62").Value = -10049.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 106380.5
$ws.Range("I65").Value = 203959.8
$ws.Range("J65").Value = 8801.200000000001
$ws.Range("K65").Value = 1019799
$ws.Range("L65").Value = 44006
$ws.Range("M65").Value = -1016679
$ws.Range("N65").Value = -50246

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 931.25
$ws.Range("I107").Value = 874.375
$ws.Range("J107").Value = 1045
$ws.Range("K107").Value = 874.375
$ws.Range("L107").Value = 1045
$ws.Range("M107").Value = 1045.625
$ws.Range("N107").Value = -4885

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1372655.4
$ws.Range("I129").Value = 281.9091
$ws.Range("J129").Value = 2316162.2
$ws.Range("K129").Value = 845.7273
$ws.Range("L129").Value = 6948486.600000001
$ws.Range("M129").Value = 4154.2727
$ws.Range("N129").Value = -6958486.600000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1530.6
$ws.Range("I132").Value = 1505.5209
$ws.Range("J132").Value = 1702.5714
$ws.Range("K132").Value = 4516.5627
$ws.Range("L132").Value = 5107.7142
$ws.Range("M132").Value = -1986.5627
$ws.Range("N132").Value = -10167.7142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 820.3077
$ws.Range("I137").Value = 742.1818
$ws.Range("K137").Value = 2226.5454
$ws.Range("M137").Value = 323.4546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3806.09
$ws.Range("I138").Value = 1910.875
$ws.Range("J138").Value = 3970.8914
$ws.Range("K138").Value = 5732.625
$ws.Range("L138").Value = 11912.6742
$ws.Range("M138").Value = -592.625
$ws.Range("N138").Value = -22192.6742

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2528.3845
$ws.Range("I141").Value = 2651.7273
$ws.Range("J141").Value = 1850
$ws.Range("K141").Value = 7955.1819
$ws.Range("L141").Value = 5550
$ws.Range("M141").Value = -2775.1819
$ws.Range("N141").Value = -15910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5641.67
$ws.Range("I32").Value = 3905.307
$ws.Range("J32").Value = 18375
$ws.Range("K32").Value = 3905.307
$ws.Range("L32").Value = 18375
$ws.Range("M32").Value = -3618.307
$ws.Range("N32").Value = -18949

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1005.6923
$ws.Range("I45").Value = 970
$ws.Range("J45").Value = 1062.8
$ws.Range("K45").Value = 970
$ws.Range("L45").Value = 1062.8
$ws.Range("M45").Value = -593
$ws.Range("N45").Value = -1816.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1196.7407
$ws.Range("I132").Value = 985.875
$ws.Range("J132").Value = 2883.6667
$ws.Range("K132").Value = 2957.625
$ws.Range("L132").Value = 8651.000100000001
$ws.Range("M132").Value = -427.625
$ws.Range("N132").Value = -13711.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 1564
$ws.Range("I97").Value = 1564
$ws.Range("K97").Value = 1564
$ws.Range("M97").Value = -573

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 10666.667
$ws.Range("J97").Value = 10666.667
$ws.Range("L97").Value = 10666.667
$ws.Range("N97").Value = -12648.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1161.6046
$ws.Range("I132").Value = 921.9487
$ws.Range("J132").Value = 3498.25
$ws.Range("K132").Value = 2765.8461
$ws.Range("L132").Value = 10494.75
$ws.Range("M132").Value = -235.8461000000002
$ws.Range("N132").Value = -15554.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 17858412
$ws.Range("I134").Value = 1109.1177
$ws.Range("K134").Value = 3327.3531
$ws.Range("M134").Value = -792.3531000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 62.875
$ws.Range("I8").Value = 62.875
$ws.Range("K8").Value = 188.625
$ws.Range("M8").Value = -49.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 60.125
$ws.Range("J12").Value = 55
$ws.Range("L12").Value = 165
$ws.Range("N12").Value = -511

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1603
$ws.Range("I68").Value = 399
$ws.Range("J68").Value = 1775
$ws.Range("K68").Value = 1197
$ws.Range("L68").Value = 5325
$ws.Range("M68").Value = -386
$ws.Range("N68").Value = -6947

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1603
$ws.Range("I71").Value = 399
$ws.Range("J71").Value = 1775
$ws.Range("K71").Value = 3591
$ws.Range("L71").Value = 15975
$ws.Range("M71").Value = 465
$ws.Range("N71").Value = -24087

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 6493760.5
$ws.Range("I97").Value = 8928834
$ws.Range("J97").Value = 233.33333
$ws.Range("K97").Value = 26786502
$ws.Range("L97").Value = 699.99999
$ws.Range("M97").Value = -26786006
$ws.Range("N97").Value = -1691.99999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 22729494
$ws.Range("J109").Value = 2986.6667
$ws.Range("L109").Value = 8960.000100000001
$ws.Range("N109").Value = -11040.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 647.413
$ws.Range("I122").Value = 539.8
$ws.Range("J122").Value = 660.53656
$ws.Range("K122").Value = 4858.2
$ws.Range("L122").Value = 5944.829040000001
$ws.Range("M122").Value = -2408.2
$ws.Range("N122").Value = -10844.82904

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 61763.47
$ws.Range("I136").Value = 126227.375
$ws.Range("K136").Value = 378682.125
$ws.Range("M136").Value = -373582.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 43004.19
$ws.Range("I137").Value = 2103.75
$ws.Range("J137").Value = 61182.168
$ws.Range("K137").Value = 6311.25
$ws.Range("L137").Value = 183546.504
$ws.Range("M137").Value = -1211.25
$ws.Range("N137").Value = -193746.504

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3456.6667
$ws.Range("I138").Value = 1920
$ws.Range("J138").Value = 4993.3335
$ws.Range("K138").Value = 5760
$ws.Range("L138").Value = 14980.0005
$ws.Range("M138").Value = -620
$ws.Range("N138").Value = -25260.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 6563.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 6563.5
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 6563.5
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -7139.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1571.1428
$ws.Range("I102").Value = 999.6
$ws.Range("K102").Value = 999.6
$ws.Range("M102").Value = 622.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2161.75
$ws.Range("I7").Value = 2158.8
$ws.Range("J7").Value = 2166.6667
$ws.Range("K7").Value = 2158.8
$ws.Range("L7").Value = 2166.6667
$ws.Range("M7").Value = -2046.8
$ws.Range("N7").Value = -2390.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 375781.12
$ws.Range("I40").Value = 482542.38
$ws.Range("J40").Value = 2116.6667
$ws.Range("K40").Value = 482542.38
$ws.Range("L40").Value = 2116.6667
$ws.Range("M40").Value = -482406.38
$ws.Range("N40").Value = -2388.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3591.25
$ws.Range("I93").Value = 3789.2307
$ws.Range("J93").Value = 2733.3333
$ws.Range("K93").Value = 3789.2307
$ws.Range("L93").Value = 2733.3333
$ws.Range("M93").Value = -2541.2307
$ws.Range("N93").Value = -5229.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2161.75
$ws.Range("I126").Value = 2158.8
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 6476.400000000001
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -4006.400000000001
$ws.Range("N126").Value = -11440.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4418.278
$ws.Range("I136").Value = 2439.3125
$ws.Range("K136").Value = 7317.9375
$ws.Range("M136").Value = -4767.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 18999.75
$ws.Range("J63").Value = 18999.75
$ws.Range("L63").Value = 18999.75
$ws.Range("N63").Value = -20247.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 14249.9
$ws.Range("J64").Value = 14249.9
$ws.Range("L64").Value = 14249.9
$ws.Range("N64").Value = -14745.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 18999.75
$ws.Range("J66").Value = 18999.75
$ws.Range("L66").Value = 56999.25
$ws.Range("N66").Value = -63239.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 14249.9
$ws.Range("J67").Value = 14249.9
$ws.Range("L67").Value = 14249.9
$ws.Range("N67").Value = -15965.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 706.2
$ws.Range("I113").Value = 528.8889
$ws.Range("J113").Value = 972.1667
$ws.Range("K113").Value = 1586.6667
$ws.Range("L113").Value = 2916.5001
$ws.Range("M113").Value = 583.3332999999998
$ws.Range("N113").Value = -7256.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 134226
$ws.Range("J139").Value = 134226
$ws.Range("L139").Value = 134226
$ws.Range("N139").Value = -144506
